$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66
$prevRow = 65

# Copy the formatting (style) of the date cell in the previous row so the
# new row's A cell gets the same date number format (style index) instead
# of Excel allocating a brand new style entry.
$ws.Cells.Item($prevRow, 1).Copy($ws.Cells.Item($row, 1))

$ws.Cells.Item($row, 1).Value = 45496
$ws.Cells.Item($row, 2).Value = 703.7508593724
$ws.Cells.Item($row, 3).Value = 241.4645322385
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 126.5683584
$ws.Cells.Item($row, 9).Value = 253.7897401611
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 20.95053483526
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 132.4973232016
$ws.Cells.Item($row, 15).Value = 59.106707803
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0.0000030048
$ws.Cells.Item($row, 18).Value = 0
$ws.Cells.Item($row, 19).Value = 0
$ws.Cells.Item($row, 20).Value = 0
$ws.Cells.Item($row, 21).Value = 347.7110712830478
$ws.Cells.Item($row, 23).Value = 0
$ws.Cells.Item($row, 24).Value = 0
$ws.Cells.Item($row, 25).Value = 0
$ws.Cells.Item($row, 26).Value = 262.05574214413
